# Fruta / hortaliza, semanal
# Insert a new weekly record as row 10, pushing existing rows 10-44 down to 11-45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (shifts old rows 10-44 down to 11-45,
# carrying their values/formatting along, as Excel does natively).
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with this week's data for Guayaba / Vega Modelo de Temuco.
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 45148
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100108
$ws.Range("H10").Value = "Tropicales y subtropicales"
$ws.Range("I10").Value = 100108001
$ws.Range("J10").Value = "Guayaba"
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 280
$ws.Range("N10").Value = 2750
$ws.Range("O10").Value = 2750
$ws.Range("P10").Value = 2750
$ws.Range("Q10").Value = "$/kilo"
$ws.Range("R10").Value = "Región de Arica y Parinacota"
$ws.Range("S10").Value = 2750
$ws.Range("T10").Value = 1
